$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from H1 (IP) onto the two new header cells so the
# existing bold/border/centered cell style (style index 1) is reused
# instead of creating a new style entry.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New header cells
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells (row 2), plain numeric values, no special style
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 5
